$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 0.1883561643835616
    "C2" = 0.5616438356164384
    "J2" = 0.03424657534246575
    "P2" = 0.1095890410958904
    "S2" = 0.1061643835616438
    "B3" = 0.01129943502824859
    "C3" = 0.06779661016949153
    "J3" = 0.05649717514124294
    "P3" = 0.7062146892655368
    "S3" = 0.1581920903954802
    "J4" = 0.08108108108108109
    "P4" = 0.5405405405405406
    "S4" = 0.3783783783783784
    "B6" = 0.06880733944954129
    "D6" = 0.02752293577981652
    "E6" = 0.004587155963302753
    "F6" = 0.07798165137614679
    "J6" = 0.2568807339449541
    "O6" = 0.01376146788990826
    "Q6" = 0.1651376146788991
    "R6" = 0.06422018348623854
    "S6" = 0.3211009174311927
    "B7" = 0.0989010989010989
    "F7" = 0.03846153846153846
    "J7" = 0.08791208791208792
    "O7" = 0.02197802197802198
    "Q7" = 0.2417582417582418
    "R7" = 0.09340659340659341
    "S7" = 0.4175824175824176
    "B8" = 0.1130820399113082
    "D8" = 0.01330376940133038
    "F8" = 0.07538802660753881
    "J8" = 0.07317073170731707
    "O8" = 0.01552106430155211
    "Q8" = 0.2084257206208426
    "R8" = 0.1019955654101996
    "S8" = 0.3991130820399113
    "B9" = 0.1204819277108434
    "D9" = 0.02409638554216868
    "F9" = 0.06024096385542169
    "J9" = 0.08433734939759036
    "O9" = 0.01204819277108434
    "Q9" = 0.144578313253012
    "R9" = 0.1144578313253012
    "S9" = 0.4397590361445783
    "B10" = 0.09937888198757763
    "D10" = 0.01863354037267081
    "E10" = 0.001552795031055901
    "F10" = 0.05279503105590062
    "J10" = 0.09239130434782608
    "O10" = 0.01319875776397516
    "Q10" = 0.2243788819875776
    "R10" = 0.1141304347826087
    "S10" = 0.3835403726708074
    "F11" = 0.003367003367003367
    "G11" = 0.1346801346801347
    "J11" = 0.09764309764309764
    "K11" = 0.2087542087542087
    "L11" = 0.531986531986532
    "S11" = 0.02356902356902357
    "G12" = 0.7283950617283951
    "J12" = 0.1975308641975309
    "K12" = 0.01851851851851852
    "L12" = 0.03703703703703703
    "S12" = 0.01851851851851852
    "G13" = 0.6808510638297872
    "J13" = 0.2978723404255319
    "S13" = 0.02127659574468085
    "J14" = 1
    "F15" = 0.02916666666666667
    "H15" = 0.1708333333333333
    "I15" = 0.0375
    "J15" = 0.3916666666666667
    "K15" = 0.04583333333333333
    "M15" = 0.008333333333333333
    "N15" = 0.004166666666666667
    "O15" = 0.08333333333333333
    "S15" = 0.2291666666666667
    "F16" = 0.03508771929824561
    "H16" = 0.1637426900584795
    "I16" = 0.07017543859649122
    "J16" = 0.3742690058479532
    "K16" = 0.1052631578947368
    "M16" = 0.02923976608187134
    "O16" = 0.1052631578947368
    "S16" = 0.1169590643274854
    "F17" = 0.03099173553719008
    "H17" = 0.1528925619834711
    "I17" = 0.08264462809917356
    "J17" = 0.4607438016528926
    "K17" = 0.08884297520661157
    "M17" = 0.02272727272727273
    "O17" = 0.0640495867768595
    "S17" = 0.09710743801652892
    "F18" = 0.04958677685950413
    "H18" = 0.1487603305785124
    "I18" = 0.07024793388429752
    "J18" = 0.4338842975206612
    "K18" = 0.09917355371900827
    "M18" = 0.02066115702479339
    "O18" = 0.07851239669421488
    "S18" = 0.09917355371900827
    "F19" = 0.01116427432216906
    "H19" = 0.2169059011164274
    "I19" = 0.07177033492822966
    "J19" = 0.3827751196172249
    "K19" = 0.1044657097288676
    "M19" = 0.0215311004784689
    "N19" = 0.0007974481658692185
    "O19" = 0.07336523125996811
    "S19" = 0.1172248803827751
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
